# Add new columns I (I0) and J (IF) to the worksheet.
# I1 = "I0" header, J1 = "IF" header, style = same bold/border style as other headers (style index 1).
# For each data row 2..37: I = 1 (constant), J = copy of the value in column H for that row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Headers
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy the style used by the other header cells (e.g. H1) onto the new headers
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats

$lastRow = 37
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 9).Value = 1                                # column I
    $ws.Cells.Item($r, 10).Value = $ws.Cells.Item($r, 8).Value2     # column J = column H
}
